$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the corrected, longer effect names
$ws.Columns.Item(1).ColumnWidth = 15.2

# Fix header casing and correct truncated/duplicated effect + level labels
$ws.Range("A1").Value = "effect"
$ws.Range("C2").Value = -1.5800207433195559
$ws.Range("E2").Value = 0.20597082564191455
$ws.Range("C3").Value = -2.021918237420485
$ws.Range("E3").Value = 0.13240124430882474
$ws.Range("C4").Value = -1.640805194535592
$ws.Range("E4").Value = 0.19382391348629088
$ws.Range("C5").Value = -1.590469133256529
$ws.Range("E5").Value = 0.20382996588496943
$ws.Range("C6").Value = -1.7434309845542093
$ws.Range("E6").Value = 0.17491922473584873
$ws.Range("C7").Value = -1.6409825400950255
$ws.Range("E7").Value = 0.1937895427237632
$ws.Range("C8").Value = -1.6423768566307746
$ws.Range("E8").Value = 0.1935195270472947
$ws.Range("C9").Value = -1.8335738797956798
$ws.Range("E9").Value = 0.15984129219753163
$ws.Range("C10").Value = -1.6912554073416481
$ws.Range("E10").Value = 0.18428802217283818
$ws.Range("C11").Value = -1.4525587389981984
$ws.Range("E11").Value = 0.2339708511781899
$ws.Range("C12").Value = -1.8989002235008299
$ws.Range("E12").Value = 0.1497332017604844
$ws.Range("C13").Value = -1.6233057754875642
$ws.Range("E13").Value = 0.1972455705580887
$ws.Range("C14").Value = -1.9070463703536606
$ws.Range("E14").Value = 0.14851840777308142
$ws.Range("C15").Value = -1.618082761034647
$ws.Range("E15").Value = 0.19827848213178556
$ws.Range("C16").Value = -1.5226079125780527
$ws.Range("E16").Value = 0.2181422485776558
$ws.Range("C17").Value = -1.779968807284091
$ws.Range("E17").Value = 0.16864340763245675
$ws.Range("C19").Value = -2.1932529746488307
$ws.Range("E19").Value = 0.1115532777725349
$ws.Range("B20").Value = "20-24"
$ws.Range("C20").Value = -2.1526195644732598
$ws.Range("E20").Value = 0.11617941932736205
$ws.Range("B21").Value = "25-29"
$ws.Range("C21").Value = -1.9802940941618412
$ws.Range("E21").Value = 0.13802863792458364
$ws.Range("B22").Value = "30-34"
$ws.Range("C22").Value = -1.906059635682643
$ws.Range("E22").Value = 0.148665028361228
$ws.Range("B23").Value = "35-39"
$ws.Range("C23").Value = -1.8967331615072418
$ws.Range("E23").Value = 0.15005803473066753
$ws.Range("B24").Value = "40-44"
$ws.Range("C24").Value = -1.7590877876552038
$ws.Range("E24").Value = 0.17220187687542304
$ws.Range("B25").Value = "45-49"
$ws.Range("C25").Value = -1.6862771204385452
$ws.Range("E25").Value = 0.18520774825067424
$ws.Range("B26").Value = "50-54"
$ws.Range("C26").Value = -1.7077692497635355
$ws.Range("E26").Value = 0.18126970937797685
$ws.Range("B27").Value = "55-59"
$ws.Range("C27").Value = -1.6580254868967692
$ws.Range("E27").Value = 0.19051478290084642
$ws.Range("B28").Value = "60-64"
$ws.Range("C28").Value = -1.5047667679439711
$ws.Range("E28").Value = 0.22206908142089196
$ws.Range("B29").Value = "65-69"
$ws.Range("C29").Value = -1.555339977396083
$ws.Range("E29").Value = 0.21111759522832077
$ws.Range("B30").Value = "70-74"
$ws.Range("C30").Value = -1.488395662918704
$ws.Range("E30").Value = 0.22573451944605594
$ws.Range("B31").Value = "75-79"
$ws.Range("C31").Value = -1.5454726123548632
$ws.Range("E31").Value = 0.2132110812163135
$ws.Range("B32").Value = "80+"
$ws.Range("C32").Value = -1.5433303973344006
$ws.Range("E32").Value = 0.21366831476845885
$ws.Range("A33").Value = "birth_control"
$ws.Range("B33").Value = "Missing"
$ws.Range("C33").Value = -1.7062138838410743
$ws.Range("E33").Value = 0.1815518694809396
$ws.Range("A34").Value = "birth_control"
$ws.Range("B34").Value = "N/A (Male)"
$ws.Range("C34").Value = -1.9070463703536606
$ws.Range("E34").Value = 0.14851840777308142
$ws.Range("A35").Value = "birth_control"
$ws.Range("B35").Value = "No"
$ws.Range("C35").Value = -1.65765752724551
$ws.Range("E35").Value = 0.19058489755281266
$ws.Range("A36").Value = "birth_control"
$ws.Range("B36").Value = "Yes"
$ws.Range("A37").Value = "cotinine_cat"
$ws.Range("B37").Value = "3+ ng/mL"
$ws.Range("A38").Value = "cotinine_cat"
$ws.Range("B38").Value = "<3 ng/mL"
$ws.Range("C38").Value = -1.7757671779129143
$ws.Range("E38").Value = 0.16935347540341225
$ws.Range("B39").Value = "Missing"
$ws.Range("C39").Value = -1.6831835347816113
$ws.Range("E39").Value = 0.18578159144309736
$ws.Range("B40").Value = "N/A (Male)"
$ws.Range("C40").Value = -1.9070463703536606
$ws.Range("E40").Value = 0.14851840777308142
$ws.Range("B41").Value = "No"
$ws.Range("C41").Value = -1.6352298292339091
$ws.Range("E41").Value = 0.19490757069347966
$ws.Range("B42").Value = "Yes"
$ws.Range("C42").Value = -1.2097903091532887
$ws.Range("E42").Value = 0.29825981522626094
$ws.Range("B43").Value = "No"
$ws.Range("C43").Value = -2.1174359498482964
$ws.Range("E43").Value = 0.12033979052935463
$ws.Range("B44").Value = "Yes"
$ws.Range("C44").Value = -1.077320662872732
$ws.Range("E44").Value = 0.3405066365866798
$ws.Range("A45").Value = "poor_sleep"
$ws.Range("B45").Value = "No"
$ws.Range("C45").Value = -1.807395504655211
$ws.Range("E45").Value = 0.16408092878777306
$ws.Range("A46").Value = "poor_sleep"
$ws.Range("B46").Value = "Yes"
$ws.Range("C46").Value = -1.7172157021781318
$ws.Range("E46").Value = 0.17956541613026022
$ws.Range("A47").Value = "short_sleep"
$ws.Range("B47").Value = "No"
$ws.Range("C47").Value = -1.7897044289554993
$ws.Range("E47").Value = 0.16700952555250514
$ws.Range("A48").Value = "short_sleep"
$ws.Range("B48").Value = "Yes"
$ws.Range("A49").Value = "sleep_med"
$ws.Range("B49").Value = "Missing"
$ws.Range("C49").Value = -1.3786286147946516
$ws.Range("E49").Value = 0.25192380084514926
$ws.Range("A50").Value = "sleep_med"
$ws.Range("B50").Value = "No"
$ws.Range("C50").Value = -1.7856349907934523
$ws.Range("E50").Value = 0.16769054523312715
$ws.Range("A51").Value = "sleep_med"
$ws.Range("B51").Value = "Yes"
$ws.Range("C51").Value = -1.6597536110576874
$ws.Range("E51").Value = 0.19018583401562628
